# Feature : Add additional cards in hash
#
# Row 7 "Aggro Rally Sword" gains a second hash entry (Anti Identifier /
# Hash Code 3), and its existing Hash Code 1 value is refreshed.
# Row 15 "Jerva/Evo Dragon" and Row 25 "Sofina/Amulet Haven" archetype
# names drop the slash in favor of a space.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Aggro Rally Sword - update Hash Code 1, and add the new
# Anti Identifier / Hash Code 3 pair (previously "None"/"None").
$ws.Range("D7").Value = "6_AQS"
$ws.Range("G7").Value = "Diamond Paladin"
$ws.Range("H7").Value = "6_FIy"

# Row 15: rename archetype "Jerva/Evo Dragon" -> "Jerva Evo Dragon"
$ws.Range("A15").Value = "Jerva Evo Dragon"

# Row 25: rename archetype "Sofina/Amulet Haven" -> "Sofina Amulet Haven"
$ws.Range("A25").Value = "Sofina Amulet Haven"

# Restore the view/selection state recorded for the sheet after the edit.
$ws.Range("A26").Select()
